# Auto-generated Excel COM-interop script
# Updates static market-price/profit figures (columns H-N) on several
# sheets of the Leves workbook, as produced by the scheduled pricing runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 8
$ws.Range("H8").Value = 119.77778
$ws.Range("I8").Value = 97
$ws.Range("J8").Value = 199.5
$ws.Range("K8").Value = 291
$ws.Range("L8").Value = 598.5
$ws.Range("M8").Value = -152
$ws.Range("N8").Value = -876.5

# Row 80
$ws.Range("H80").Value = 493.55554
$ws.Range("I80").Value = 168.4
$ws.Range("J80").Value = 900
$ws.Range("K80").Value = 505.2
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = 492.8
$ws.Range("N80").Value = -4696

# Row 83
$ws.Range("H83").Value = 493.55554
$ws.Range("I83").Value = 168.4
$ws.Range("J83").Value = 900
$ws.Range("K83").Value = 1515.6
$ws.Range("L83").Value = 8100
$ws.Range("M83").Value = 3476.4
$ws.Range("N83").Value = -18084

# Row 95
$ws.Range("H95").Value = 37099.25
$ws.Range("J95").Value = 37099.25
$ws.Range("L95").Value = 37099.25
$ws.Range("N95").Value = -42591.25

# Row 137
$ws.Range("H137").Value = 2553.4
$ws.Range("I137").Value = 1717
$ws.Range("J137").Value = 3575.6667
$ws.Range("K137").Value = 5151
$ws.Range("L137").Value = 10727.0001
$ws.Range("M137").Value = -2601
$ws.Range("N137").Value = -15827.0001

# Row 138
$ws.Range("H138").Value = 5659.6045
$ws.Range("I138").Value = 1133.3334
$ws.Range("J138").Value = 5999.075
$ws.Range("K138").Value = 3400.0002
$ws.Range("L138").Value = 17997.225
$ws.Range("M138").Value = 1739.9998
$ws.Range("N138").Value = -28277.225

$ws = $wb.Worksheets.Item("ARM")

# Row 102
$ws.Range("H102").Value = 10067.5
$ws.Range("I102").Value = 15175
$ws.Range("K102").Value = 15175
$ws.Range("M102").Value = -13553

# Row 122
$ws.Range("H122").Value = 1906.2
$ws.Range("I122").Value = 1863
$ws.Range("J122").Value = 1986.4286
$ws.Range("K122").Value = 5589
$ws.Range("L122").Value = 5959.2858
$ws.Range("M122").Value = -3139
$ws.Range("N122").Value = -10859.2858

$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 2707.1428
$ws.Range("I20").Value = 3129.3
$ws.Range("K20").Value = 3129.3
$ws.Range("M20").Value = -2882.3

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 134
$ws.Range("H134").Value = 2369.7273
$ws.Range("I134").Value = 2260.0527
$ws.Range("J134").Value = 3064.3333
$ws.Range("K134").Value = 6780.158100000001
$ws.Range("L134").Value = 9192.999899999999
$ws.Range("M134").Value = -4245.158100000001
$ws.Range("N134").Value = -14262.9999

# Row 139
$ws.Range("H139").Value = 80709
$ws.Range("I139").Value = 80709
$ws.Range("K139").Value = 80709
$ws.Range("M139").Value = -75569

$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 73999.75
$ws.Range("I16").Value = 76500
$ws.Range("J16").Value = 71499.5
$ws.Range("K16").Value = 76500
$ws.Range("L16").Value = 71499.5
$ws.Range("M16").Value = -76213
$ws.Range("N16").Value = -72073.5

# Row 31
$ws.Range("H31").Value = 6670.9473
$ws.Range("I31").Value = 8241.637000000001
$ws.Range("K31").Value = 8241.637000000001
$ws.Range("M31").Value = -7946.637000000001

# Row 34
$ws.Range("H34").Value = 6670.9473
$ws.Range("I34").Value = 8241.637000000001
$ws.Range("K34").Value = 8241.637000000001
$ws.Range("M34").Value = -8039.637000000001

# Row 59
$ws.Range("H59").Value = 85998.60000000001
$ws.Range("J59").Value = 99998.25
$ws.Range("L59").Value = 99998.25
$ws.Range("N59").Value = -102288.25

# Row 105
$ws.Range("H105").Value = 1640.2
$ws.Range("I105").Value = 1550.25
$ws.Range("K105").Value = 1550.25
$ws.Range("M105").Value = 196.75

# Row 113
$ws.Range("H113").Value = 73999.75
$ws.Range("I113").Value = 76500
$ws.Range("J113").Value = 71499.5
$ws.Range("K113").Value = 76500
$ws.Range("L113").Value = 71499.5
$ws.Range("M113").Value = -74330
$ws.Range("N113").Value = -75839.5

# Row 134
$ws.Range("H134").Value = 2843
$ws.Range("I134").Value = 2843
$ws.Range("K134").Value = 8529
$ws.Range("M134").Value = -5994

$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 3829.8333
$ws.Range("I3").Value = 3829.8333
$ws.Range("K3").Value = 11489.4999
$ws.Range("M3").Value = -11377.4999

# Row 68
$ws.Range("H68").Value = 1856.7273
$ws.Range("I68").Value = 1300.6666
$ws.Range("J68").Value = 2065.25
$ws.Range("K68").Value = 3901.9998
$ws.Range("L68").Value = 6195.75
$ws.Range("M68").Value = -3090.9998
$ws.Range("N68").Value = -7817.75

# Row 71
$ws.Range("H71").Value = 1856.7273
$ws.Range("I71").Value = 1300.6666
$ws.Range("J71").Value = 2065.25
$ws.Range("K71").Value = 11705.9994
$ws.Range("L71").Value = 18587.25
$ws.Range("M71").Value = -7649.999400000001
$ws.Range("N71").Value = -26699.25

# Row 98
$ws.Range("H98").Value = 740.41174
$ws.Range("J98").Value = 380.6
$ws.Range("L98").Value = 1141.8
$ws.Range("N98").Value = -4137.8

# Row 107
$ws.Range("H107").Value = 764.04
$ws.Range("I107").Value = 564.55554
$ws.Range("J107").Value = 876.25
$ws.Range("K107").Value = 1693.66662
$ws.Range("L107").Value = 2628.75
$ws.Range("M107").Value = 226.33338
$ws.Range("N107").Value = -6468.75

# Row 139
$ws.Range("H139").Value = 3535.4167
$ws.Range("I139").Value = 2677.3684
$ws.Range("J139").Value = 6796
$ws.Range("K139").Value = 8032.1052
$ws.Range("L139").Value = 20388
$ws.Range("M139").Value = -2892.1052
$ws.Range("N139").Value = -30668

$ws = $wb.Worksheets.Item("LTW")

# Row 68
$ws.Range("H68").Value = 5500
$ws.Range("I68").Value = 5500
$ws.Range("K68").Value = 5500
$ws.Range("M68").Value = -4751

# Row 71
$ws.Range("H71").Value = 5500
$ws.Range("I71").Value = 5500
$ws.Range("K71").Value = 27500
$ws.Range("M71").Value = -23756

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

# Row 100
$ws.Range("H100").Value = 9999.5
$ws.Range("I100").Value = 9999.5
$ws.Range("K100").Value = 9999.5
$ws.Range("M100").Value = -9458.5

$ws = $wb.Worksheets.Item("WVR")

# Row 126
$ws.Range("H126").Value = 3894.182
$ws.Range("I126").Value = 4182.579
$ws.Range("J126").Value = 2067.6667
$ws.Range("K126").Value = 12547.737
$ws.Range("L126").Value = 6203.000100000001
$ws.Range("M126").Value = -10077.737
$ws.Range("N126").Value = -11143.0001
